$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# The sheet holds a "missing items" table. Two new product rows are
# being added to the daily report:
#   - "LARYPRO 20 LOZENGES"  -> inserted right before "MIDODRINE 2.5MG 20 TAB"
#   - "حبايه"                 -> inserted right before "زيت جونسون صغير"
# The running total (originally on row 12) and the generated-at footer
# (originally row 13) shift down by the two inserted rows, and the
# timestamp in the footer is refreshed.
# ----------------------------------------------------------------------

# Helper: write a value into a cell while preserving its existing
# number format. Excel otherwise happily re-interprets numeric looking
# text such as "22.0000" or even "1" as a real number and drops the
# shared-string storage the template relies on. Cells that are already
# formatted as plain text ("@") are left alone, since re-applying "@"
# on some of them snaps the cell to a different (but visually
# equivalent) style id and loses the original reading-order flag.
function Set-TextValue($range, $value) {
    $fmt = $range.NumberFormat
    if ($fmt -ne "@") {
        $range.NumberFormat = "@"
        $range.Value2 = $value
        $range.NumberFormat = $fmt
    } else {
        $range.Value2 = $value
    }
}

# Insert a blank row just above the MIDODRINE row (row 9).
$ws.Rows.Item(9).Insert()
# Insert a blank row just above the زيت جونسون صغير row (now row 11,
# since the previous insert pushed it down by one).
$ws.Rows.Item(11).Insert()

# ---- Row 9: copy formatting (incl. merges) from row 10 (MIDODRINE) ----
$ws.Range("A10:Q10").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A9").Value2 = 3
$ws.Range("C9").Value2 = "LARYPRO 20 LOZENGES"
Set-TextValue $ws.Range("H9") "1:0"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "44.00"
Set-TextValue $ws.Range("P9") "22.0000"
Set-TextValue $ws.Range("Q9") "0:1"

# Renumber MIDODRINE (now row 10) to keep the "م" sequence intact.
$ws.Range("A10").Value2 = 4

# ---- Row 11: copy formatting (incl. merges) from row 12 (زيت جونسون) ----
$ws.Range("A12:Q12").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

$ws.Range("A11").Value2 = 5
$ws.Range("C11").Value2 = "حبايه"
Set-TextValue $ws.Range("H11") "0:0"
Set-TextValue $ws.Range("L11") "0"
Set-TextValue $ws.Range("N11") "3.00"
Set-TextValue $ws.Range("P11") "6.0000"
Set-TextValue $ws.Range("Q11") "2:0"

# Renumber the rows that followed (زيت جونسون صغير, قطن 50جم).
$ws.Range("A12").Value2 = 6
$ws.Range("A13").Value2 = 7

# Row heights: rows 9-13 alternate 25.5 / 24.75 per the template, with
# row 13 (now قطن 50جم, the last product row) set to 24.75 as well.
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75

# ---- Update the running total (now row 14) ----
$ws.Range("P14").Value2 = 269

# ---- Refresh the generated-at timestamp in the footer (now row 15) ----
$ws.Range("A15").Value2 = "Sunday, 22 June, 2025 10:03 AM"
